$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the manualStatus (column I) values in rows 5-7 from numeric 64 to text "[64]"
$ws.Range("I5").Value = "[64]"
$ws.Range("I6").Value = "[64]"
$ws.Range("I7").Value = "[64]"

# Widen column F (fastqFileName) to fit content
$ws.Columns.Item(6).ColumnWidth = 48.75

# Row heights for rows 6 and 7 settle slightly shorter after the edit
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8

# Move the active cell selection to I7
$ws.Range("I7").Select()
